# Edit the "Test Suite" sheet: set Runmode column (C) to "N" for all
# data rows (3-37), except row 32 ("Dashboard Suite") which stays "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

for ($r = 3; $r -le 37; $r++) {
    if ($r -eq 32) {
        continue
    }
    $ws.Cells.Item($r, 3).Value = "N"
}

# Update the active selection to mirror the author's final view state:
# C32 becomes the active/selected cell (also nudges the frozen pane's
# scroll position forward, matching topLeftCell moving from A18 to A12).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("C32").Select()
